$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 263.42
$ws.Range("I15").Value = 263.42
$ws.Range("K15").Value = 790.26
$ws.Range("M15").Value = -621.26
# Row 92
$ws.Range("H92").Value = 22224516
$ws.Range("I92").Value = 30304678
$ws.Range("J92").Value = 4072.5
$ws.Range("K92").Value = 30304678
$ws.Range("L92").Value = 4072.5
$ws.Range("M92").Value = -30303430
$ws.Range("N92").Value = -6568.5
# Row 96
$ws.Range("H96").Value = 2147.5
$ws.Range("I96").Value = 1095
$ws.Range("J96").Value = 3200
$ws.Range("K96").Value = 3285
$ws.Range("L96").Value = 9600
$ws.Range("M96").Value = -1912
$ws.Range("N96").Value = -12346
# Row 98
$ws.Range("H98").Value = 963.381
$ws.Range("I98").Value = 963.381
$ws.Range("K98").Value = 963.381
$ws.Range("M98").Value = 534.619
# Row 111
$ws.Range("H111").Value = 4614.5
$ws.Range("I111").Value = 3886
$ws.Range("J111").Value = 6800
$ws.Range("K111").Value = 11658
$ws.Range("L111").Value = 20400
$ws.Range("M111").Value = -8591
$ws.Range("N111").Value = -26534
# Row 112
$ws.Range("H112").Value = 4870.8
$ws.Range("J112").Value = 5132.766
$ws.Range("L112").Value = 15398.298
$ws.Range("N112").Value = -17614.298
# Row 121
$ws.Range("H121").Value = 1240.8064
$ws.Range("J121").Value = 1240.8064
$ws.Range("L121").Value = 3722.4192
$ws.Range("N121").Value = -7216.4192
# Row 122
$ws.Range("H122").Value = 963.381
$ws.Range("I122").Value = 963.381
$ws.Range("K122").Value = 2890.143
$ws.Range("M122").Value = -440.143
# Row 140
$ws.Range("H140").Value = 76790.25
$ws.Range("J140").Value = 76790.25
$ws.Range("L140").Value = 76790.25
$ws.Range("N140").Value = -87150.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2800
$ws.Range("I2").Value = 3005
$ws.Range("J2").Value = 2697.5
$ws.Range("K2").Value = 3005
$ws.Range("L2").Value = 2697.5
$ws.Range("M2").Value = -2892
$ws.Range("N2").Value = -2923.5
# Row 32
$ws.Range("H32").Value = 66700470
$ws.Range("I32").Value = 71462720
$ws.Range("J32").Value = 29000
$ws.Range("K32").Value = 71462720
$ws.Range("L32").Value = 29000
$ws.Range("M32").Value = -71462433
$ws.Range("N32").Value = -29574
# Row 45
$ws.Range("H45").Value = 2157.5264
$ws.Range("I45").Value = 1422.3
$ws.Range("J45").Value = 2974.4443
$ws.Range("K45").Value = 1422.3
$ws.Range("L45").Value = 2974.4443
$ws.Range("M45").Value = -1045.3
$ws.Range("N45").Value = -3728.4443
# Row 74
$ws.Range("H74").Value = 23811446
$ws.Range("I74").Value = 1431.2
$ws.Range("J74").Value = 83336480
$ws.Range("K74").Value = 1431.2
$ws.Range("L74").Value = 83336480
$ws.Range("M74").Value = -557.2
$ws.Range("N74").Value = -83338228
# Row 77
$ws.Range("H77").Value = 23811446
$ws.Range("I77").Value = 1431.2
$ws.Range("J77").Value = 83336480
$ws.Range("K77").Value = 7156
$ws.Range("L77").Value = 416682400
$ws.Range("M77").Value = -2788
$ws.Range("N77").Value = -416691136
# Row 116
$ws.Range("H116").Value = 2800
$ws.Range("I116").Value = 3005
$ws.Range("J116").Value = 2697.5
$ws.Range("K116").Value = 3005
$ws.Range("L116").Value = 2697.5
$ws.Range("M116").Value = -711
$ws.Range("N116").Value = -7285.5
# Row 122
$ws.Range("H122").Value = 76094.375
$ws.Range("I122").Value = 93344.766
$ws.Range("J122").Value = 1342.6666
$ws.Range("K122").Value = 280034.298
$ws.Range("L122").Value = 4027.9998
$ws.Range("M122").Value = -277584.298
$ws.Range("N122").Value = -8927.9998

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2800
$ws.Range("I3").Value = 3005
$ws.Range("J3").Value = 2697.5
$ws.Range("K3").Value = 3005
$ws.Range("L3").Value = 2697.5
$ws.Range("M3").Value = -2891
$ws.Range("N3").Value = -2925.5
# Row 94
$ws.Range("H94").Value = 100522.2
$ws.Range("I94").Value = 250362.25
$ws.Range("J94").Value = 628.8333
$ws.Range("K94").Value = 250362.25
$ws.Range("L94").Value = 628.8333
$ws.Range("M94").Value = -249911.25
$ws.Range("N94").Value = -1530.8333

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 9957.812
$ws.Range("I31").Value = 1487.4286
$ws.Range("J31").Value = 11246.782
$ws.Range("K31").Value = 1487.4286
$ws.Range("L31").Value = 11246.782
$ws.Range("M31").Value = -1192.4286
$ws.Range("N31").Value = -11836.782
# Row 34
$ws.Range("H34").Value = 9957.812
$ws.Range("I34").Value = 1487.4286
$ws.Range("J34").Value = 11246.782
$ws.Range("K34").Value = 1487.4286
$ws.Range("L34").Value = 11246.782
$ws.Range("M34").Value = -1285.4286
$ws.Range("N34").Value = -11650.782
# Row 70
$ws.Range("H70").Value = 29000
$ws.Range("I70").Value = 20000
$ws.Range("J70").Value = 38000
$ws.Range("K70").Value = 20000
$ws.Range("L70").Value = 38000
$ws.Range("M70").Value = -19685
$ws.Range("N70").Value = -38630
# Row 73
$ws.Range("H73").Value = 29000
$ws.Range("I73").Value = 20000
$ws.Range("J73").Value = 38000
$ws.Range("K73").Value = 20000
$ws.Range("L73").Value = 38000
$ws.Range("M73").Value = -18908
$ws.Range("N73").Value = -40184

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1048.9734
$ws.Range("I5").Value = 801.4761999999999
$ws.Range("J5").Value = 2348.3333
$ws.Range("K5").Value = 2404.4286
$ws.Range("L5").Value = 7044.999899999999
$ws.Range("M5").Value = -2292.4286
$ws.Range("N5").Value = -7268.999899999999
# Row 14
$ws.Range("H14").Value = 41666830
$ws.Range("I14").Value = 41666830
$ws.Range("K14").Value = 125000490
$ws.Range("M14").Value = -125000317
# Row 58
$ws.Range("H58").Value = 1914.6428
$ws.Range("I58").Value = 805
$ws.Range("K58").Value = 2415
$ws.Range("M58").Value = -2287
# Row 107
$ws.Range("H107").Value = 18182724
$ws.Range("I107").Value = 260.33334
$ws.Range("J107").Value = 27028246
$ws.Range("K107").Value = 781.0000200000001
$ws.Range("L107").Value = 81084738
$ws.Range("M107").Value = 1138.99998
$ws.Range("N107").Value = -81088578
# Row 113
$ws.Range("H113").Value = 852.9178000000001
$ws.Range("I113").Value = 703.6889
$ws.Range("J113").Value = 1092.75
$ws.Range("K113").Value = 2111.0667
$ws.Range("L113").Value = 3278.25
$ws.Range("M113").Value = 58.93330000000014
$ws.Range("N113").Value = -7618.25
# Row 131
$ws.Range("H131").Value = 4091.7
$ws.Range("I131").Value = 1000
$ws.Range("J131").Value = 4254.421
$ws.Range("K131").Value = 3000
$ws.Range("L131").Value = 12763.263
$ws.Range("M131").Value = 2040
$ws.Range("N131").Value = -22843.263
# Row 132
$ws.Range("H132").Value = 2448.3394
$ws.Range("I132").Value = 2583.6924
$ws.Range("J132").Value = 2331.0334
$ws.Range("K132").Value = 23253.2316
$ws.Range("L132").Value = 20979.3006
$ws.Range("M132").Value = -20723.2316
$ws.Range("N132").Value = -26039.3006
# Row 135
$ws.Range("H135").Value = 1048.9734
$ws.Range("I135").Value = 801.4761999999999
$ws.Range("J135").Value = 2348.3333
$ws.Range("K135").Value = 7213.2858
$ws.Range("L135").Value = 21134.9997
$ws.Range("M135").Value = -4678.2858
$ws.Range("N135").Value = -26204.9997

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 670887.4399999999
$ws.Range("I3").Value = 1003796.6
$ws.Range("J3").Value = 254751
$ws.Range("K3").Value = 1003796.6
$ws.Range("L3").Value = 254751
$ws.Range("M3").Value = -1003680.6
$ws.Range("N3").Value = -254983
# Row 15
$ws.Range("H15").Value = 10105.263
$ws.Range("J15").Value = 10105.263
$ws.Range("L15").Value = 10105.263
$ws.Range("N15").Value = -10681.263
# Row 81
$ws.Range("H81").Value = 10105.263
$ws.Range("J81").Value = 10105.263
$ws.Range("L81").Value = 10105.263
$ws.Range("N81").Value = -12101.263
# Row 84
$ws.Range("H84").Value = 10105.263
$ws.Range("J84").Value = 10105.263
$ws.Range("L84").Value = 30315.789
$ws.Range("N84").Value = -40299.789
# Row 113
$ws.Range("H113").Value = 62990.89
$ws.Range("I113").Value = 86479.62
$ws.Range("K113").Value = 86479.62
$ws.Range("M113").Value = -84309.62

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 80
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
# Row 83
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 2699
$ws.Range("I122").Value = 2408.8333
$ws.Range("J122").Value = 2966.8462
$ws.Range("K122").Value = 7226.499899999999
$ws.Range("L122").Value = 8900.5386
$ws.Range("M122").Value = -4776.499899999999
$ws.Range("N122").Value = -13800.5386
